# issue #5: add legislator_id, name, date into dataframe
#
# The 股票 (stocks) worksheet (sheet4) gains three new trailing columns:
#   H = date             -> "2012-04-30" (literal text, not a date serial)
#   I = legislator_name  -> "徐少萍"
#   J = legislator_id    -> 726
#
# Every existing data row (2-33) gets the same three values; the header
# row (row 1) gets the three new column names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$reportDate = "2012-04-30"
$legislatorName = "徐少萍"
$legislatorId = 726

# Header row
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Determine the last populated data row from column A (rows 2..33 in this sheet).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    # Force the date column to stay plain text (shared string) instead of
    # being auto-converted into a date serial number.
    $cell = $ws.Cells.Item($r, 8)
    $cell.NumberFormat = "@"
    $cell.Value = $reportDate

    $ws.Cells.Item($r, 9).Value = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
